# Add a "Result" / "Pass" column (F) to the Sheet1 data table:
#   F1 = "Result" header (same style as the other headers)
#   F2 = "Pass"   value  (left unstyled, like the other appended cell)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inserting column F first makes it inherit the existing row formatting
# (style used by A1:E2) so the new header cell lines up with the rest of
# the header row.
$ws.Columns.Item(6).Insert()

$ws.Range("F1").Value = "Result"
$ws.Range("F2").Value = "Pass"

# The data cell keeps the default (unstyled) format.
$ws.Range("F2").Style = "Normal"

# Match the workbook's final selection/active cell.
$ws.Range("F2").Select()
